# Alunos.xlsx — "Add files via upload": append a test student row (26)
# to the roster: name, GRR login, e-mail (as a mailto hyperlink) and the
# GRR-number formula carried down from the row above, formatted like the
# rest of the table (wrapped text + a boxed left/right border).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new row data ---------------------------------------------------
$ws.Range("A26").Value = "ALUNO TESTE"
$ws.Range("B26").Value = "GRR20209999"
$ws.Range("C26").Value = "teste@ufpr.br"
$ws.Range("D26").Formula = "=RIGHT(B26,8)"

# ---- formatting: medium left/right border + wrapped text, per cell --
$rowCells = @("A26", "B26", "C26", "D26")
foreach ($addr in $rowCells) {
    $cell = $ws.Range($addr)
    $cell.WrapText = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = -4138
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = -4138
}

$ws.Rows.Item(26).RowHeight = 14

# ---- turn the e-mail into a mailto hyperlink -------------------------
$ws.Hyperlinks.Add($ws.Range("C26"), "mailto:teste@ufpr.br")

# ---- match the saved selection from the source session --------------
$ws.Range("C13").Select()
